$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine last used row in column A (data rows start at row 2)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $areaCell = $ws.Cells.Item($r, 2)
    $stationCell = $ws.Cells.Item($r, 3)

    if ($areaCell.Value2 -eq "ITA17") {
        $areaCell.Value = "SLO"
    }
    if ($stationCell.Value2 -eq "45BIS") {
        $stationCell.Value = "45bis"
    }
}

# Update Oth (column L) for rows 60-63 from 0 to 2
foreach ($r in 60..63) {
    $othCell = $ws.Cells.Item($r, 12)
    if ($othCell.Value2 -eq 0) {
        $othCell.Value = 2
    }
}
